# Ajustes para inserir checklist
# Rename the three "Materiais" category labels to the new checklist items.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Alternador"
$ws.Range("A2").Value = "Airbag"
$ws.Range("A3").Value = "Ala do teto"

# Reset the active selection back to the top-left cell.
[void]$ws.Range("A1").Select()
